$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.363.68"
$ws.Range("E2").Value = "  -2.24%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.643.35"
$ws.Range("E3").Value = "  -3.30%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.36"
$ws.Range("E5").Value = "  -0.83%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.19"
$ws.Range("E6").Value = "  -1.68%  "

# Row 8
$ws.Range("E8").Value = "  -0.80%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.644.01"
$ws.Range("E9").Value = "  -3.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("E10").Value = "  -0.82%  "

# Row 12
$ws.Range("E12").Value = "  -0.48%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.23"
$ws.Range("E13").Value = "  -2.12%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.99"
$ws.Range("E14").Value = "  -2.53%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.124.28"
$ws.Range("E15").Value = "  -3.34%  "

# Row 16
$ws.Range("E16").Value = "  -3.20%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.217.34"
$ws.Range("E17").Value = "  -2.25%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.645.15"
$ws.Range("E18").Value = "  -2.97%  "

# Row 19
$ws.Range("E19").Value = "  +0.53%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.89"
$ws.Range("E20").Value = "  +2.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.91"
$ws.Range("E21").Value = "  -2.63%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.40"
$ws.Range("E22").Value = "  -2.97%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  +11.42%  "

# Row 25
$ws.Range("E25").Value = "  -6.38%  "

# Row 26
$ws.Range("E26").Value = "  +0.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.88"
$ws.Range("E27").Value = "  -4.04%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.777.95"
$ws.Range("E28").Value = "  -3.44%  "

# Row 29
$ws.Range("E29").Value = "  -3.88%  "

# Row 30
$ws.Range("E30").Value = "  +0.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "554.72"
$ws.Range("E31").Value = "  -6.55%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("E32").Value = "  -2.96%  "

# Row 33
$ws.Range("E33").Value = "  -4.01%  "

# Row 34
$ws.Range("E34").Value = "  -1.84%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("E36").Value = "  -0.03%  "

# Row 37
$ws.Range("E37").Value = "  -5.29%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.66"
$ws.Range("E38").Value = "  -2.36%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.41"
$ws.Range("E39").Value = "  -2.37%  "

# Row 40
$ws.Range("E40").Value = "  -2.47%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.27"
$ws.Range("E41").Value = "  -4.32%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").Value = "  -5.13%  "

# Row 43
$ws.Range("E43").Value = "  -0.47%  "

# Row 44
$ws.Range("E44").Value = "  -5.92%  "

# Row 45
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.15"
$ws.Range("E46").Value = "  -2.12%  "

# Row 47
$ws.Range("E47").Value = "  -3.57%  "

# Row 48
$ws.Range("E48").Value = "  -1.51%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "154.13"
$ws.Range("E49").Value = "  -1.88%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.89"
$ws.Range("E50").Value = "  -2.02%  "

# Row 51
$ws.Range("E51").Value = "  -3.57%  "
